$d = $word.ActiveDocument

$replacements = @(
    @("13÷7=1, 6",  "73÷8=9, 1"),
    @("96÷4=24, 0", "81÷3=27, 0"),
    @("91÷7=13, 0", "66÷5=13, 1"),
    @("55÷2=27, 1", "99÷5=19, 4"),
    @("83÷8=10, 3", "77÷4=19, 1"),
    @("23÷7=3, 2",  "63÷3=21, 0"),
    @("66÷4=16, 2", "45÷2=22, 1"),
    @("30÷9=3, 3",  "99÷4=24, 3"),
    @("86÷4=21, 2", "40÷6=6, 4"),
    @("24÷9=2, 6",  "32÷9=3, 5"),
    @("76÷9=8, 4",  "46÷2=23, 0"),
    @("84÷6=14, 0", "38÷5=7, 3"),
    @("79÷9=8, 7",  "79÷4=19, 3"),
    @("92÷2=46, 0", "90÷2=45, 0"),
    @("83÷7=11, 6", "88÷7=12, 4"),
    @("81÷6=13, 3", "13÷5=2, 3"),
    @("90÷5=18, 0", "95÷9=10, 5"),
    @("66÷2=33, 0", "48÷3=16, 0"),
    @("38÷4=9, 2",  "96÷3=32, 0"),
    @("11÷2=5, 1",  "44÷9=4, 8"),
    @("55÷3=18, 1", "82÷8=10, 2"),
    @("94÷3=31, 1", "44÷2=22, 0"),
    @("20÷7=2, 6",  "51÷8=6, 3"),
    @("74÷8=9, 2",  "50÷9=5, 5"),
    @("97÷2=48, 1", "12÷3=4, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
